# The commit ("Modified the data ...") removes the first two data rows of
# the table on Sheet1 (old rows 2 and 3), so every later row shifts up by
# two positions and the used range shrinks from A1:D74 to A1:D72. The
# B/C/D formulas are relative, so they recompute automatically once the
# rows are gone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first two data rows (the header in row 1 is left untouched).
$ws.Rows("2:3").Delete()

# Move the view/selection to match the post-edit cursor position.
$ws.Range("I12").Select()

# Best-effort: also update the scroll position of the window so the top
# visible row lines up with where the author left the sheet scrolled to.
$win = $excel.ActiveWindow
if ($win) {
    $win.ScrollRow = 61
    $win.ScrollColumn = 1
}
